$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("s1")
$ws.Range("B5:G5").Copy()
$ws.Range("B13").PasteSpecial(-4122)
Write-Host ("pasted")
